# Add the new paper entry as described in the commit:
# "Eerste extra paper toegevoegd, zonder comments"
# (First extra paper added, without comments)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3: only column A is filled in (paper name), columns B and C
# (Interesting Content / What can it be used for?) are left empty.
$ws.Range("A3").Value = "Review of the Robustness and Applicability of Monocular Pose Estimation Systems for Relative Navigation with an Uncooperative Spacecraft"

# Column A needs to widen to fit the long text (no longer "best fit",
# now an explicit custom width). The runtime quantizes ColumnWidth to
# 1/6-character steps, so 115.3333... is the input that lands on the
# stored width closest to the target 116.109375 (i.e. 116.166666...).
$ws.Columns.Item(1).ColumnWidth = 115.3333333333
